$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Housing Dataset")

# Insert a new "STD" column into each of the two result blocks.
# Block 1 currently spans F:I (Training/CV/Test/Describe) -> insert before H (Test)
# Block 2 currently spans J:M (Training/CV/Test/Describe) -> insert before the
# (post-shift) second "Test" column, which is column M once the first insert
# has shifted the sheet right by one.
$ws.Columns("H").Insert()
$ws.Columns("M").Insert()

# New header labels for the inserted "STD" columns.
$ws.Range("H2").Value = "STD"
$ws.Range("M2").Value = "STD"

# New data points for row 3 (the "#1 / AS IS" scenario row).
$ws.Range("H3").Value = 1439
$ws.Range("M3").Value = 672

# Updated learning-curve description text (was "good learning curve = little
# variance", now "good learning curve = little to no variance").
$ws.Range("O3").Value = "good learning curve = little to no variance"

# The formerly-empty, styled B1 cell is gone in the edited workbook.
$ws.Range("B1").Clear()
